# Task1.xlsx edit: add a title row with rich-text header, bold the column
# header row, and fix the postfix/operator-stack trace for the infix
# expression a*b/(c-a)+d*e (the '/' is now flushed into the postfix output
# as soon as '+' forces it off the stack, instead of being carried all the
# way to the end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row 1 for the title, pushing everything else down by one.
$ws.Rows("1:1").Insert()

# 2. Title cell: rich text, first part bold ("Infix Expression"), second
#    part regular (": a*b/(c-a)+d*e"). Merge A1:C1 and center it.
$ws.Range("A1").Value = "Infix Expression: a*b/(c-a)+d*e"
$titleLen = "Infix Expression".Length
$ws.Range("A1").Characters(1, $titleLen).Font.Bold = $true
$ws.Range("A1:C1").Merge()
$ws.Range("A1:C1").HorizontalAlignment = -4108

# 3. Header row (now row 2) becomes bold.
$ws.Range("A2:C2").Font.Bold = $true

# 4. Fix the data rows affected by the corrected trace (rows 13-18 after
#    the insert correspond to the steps once '+' is seen and the pending
#    '/' is popped into the postfix output).
$ws.Range("A13").Value = "+"
$ws.Range("B13").Value = "ab*ca-/"
$ws.Range("C13").Value = "+"

$ws.Range("A14").Value = "d"
$ws.Range("B14").Value = "ab*ca-/d"
$ws.Range("C14").Value = "+"

$ws.Range("A15").Value = "e"
$ws.Range("B15").Value = "ab*ca-/de"
$ws.Range("C15").Value = "+"

$ws.Range("A16").Value = "*"
$ws.Range("B16").Value = "ab*ca-/de"
$ws.Range("C16").Value = "+*"
$ws.Range("C16").NumberFormat = "@"

$ws.Range("B17").Value = "ab*ca-/de*"
$ws.Range("C17").Value = "+"

$ws.Range("B18").Value = "ab*ca-/de*+"
$ws.Range("C18").ClearContents()

# 5. The old final row (now row 19, holding the stale trailing '/' pop) is
#    no longer needed.
$ws.Rows("19:19").Delete()

# 6. Misc formatting to match the saved state.
$ws.PageSetup.Orientation = 1
$ws.Range("C8").Select()
